$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product code column (A3:A6) from SP0006 to SP0003
$ws.Range("A3:A6").Value = "SP0003"

# Update quantity in C3 from 2 to 80
$ws.Cells.Item(3, 3).Value = 80

# Add new row 7 with demo variant data, copying from row 4 so that
# the text-looking numeric values ("100.3" / "222.7") stay shared strings
$ws.Range("A4").Copy($ws.Range("A7"))
$ws.Cells.Item(7, 2).Value = "Demo"
$ws.Cells.Item(7, 3).Value = 80
$ws.Range("D4").Copy($ws.Range("D7"))
$ws.Range("E4").Copy($ws.Range("E7"))
$ws.Range("F4").Copy($ws.Range("F7"))

# Widen column F
$ws.Columns.Item(6).ColumnWidth = 14.65

# Update active selection
$ws.Range("E8").Select()
